# Weekly update: insert a new price record for "Poroto granado" at row 132
# (Vega Central Mapocho de Santiago, origin "Perú"), shifting the existing
# rows 132-155 down to 133-156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 132; Excel automatically shifts the
# following rows (and the sheet dimension) down by one.
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A132").Value = 9
$ws.Range("B132").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C132").Value = "Metropolitana"
$ws.Range("D132").Value = 44504
$ws.Range("E132").Value = 13
$ws.Range("F132").Value = 100112030
$ws.Range("G132").Value = "Poroto granado"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 44
$ws.Range("K132").Value = 31000
$ws.Range("L132").Value = 33000
$ws.Range("M132").Value = 32000
$ws.Range("N132").Value = "$/malla 25 kilos"
$ws.Range("O132").Value = "Perú"
$ws.Range("P132").Value = 1280
$ws.Range("Q132").Value = 25
$ws.Range("R132").Value = "Hortaliza"
